# Update column G ("K") values for rows 2-36 to reflect the regenerated
# save_data (K computed instead of Strike#, using recalculated std/mean and
# s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4,2,0,2,0,1,2,1,2,2,2,2,4,1,3,5,3,5,2,1,2,4,5,5,1,2,3,3,6,6,8,2,4,4,3)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
